# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change: cell B11 on the "Rules" sheet is updated from "R40" to the
# text value "1" (kept as a text/string cell, same cell style as before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# Write the value through a formula that evaluates to the text string "1"
# and then convert the formula result into a plain (non-formula) value via
# copy / paste-values. This keeps the cell typed as text (matches the
# original shared-string cell type) without Excel re-interpreting the
# literal "1" as a number, and without altering the cell's existing style.
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
